$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '46.423.81'
$ws.Range('E2').Value = '  -2.02%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.616.91'
$ws.Range('E3').Value = '  +1.05%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.29%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.80'
$ws.Range('E5').Value = '  +0.23%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.28'
$ws.Range('E6').Value = '  -4.13%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.603'
$ws.Range('E7').Value = '  -1.22%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.14%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  +0.51%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.53'
$ws.Range('E10').Value = '  +0.80%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.31'
$ws.Range('E11').Value = '  -1.02%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0843'
$ws.Range('E12').Value = '  +0.80%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.15'
$ws.Range('E13').Value = '  +1.36%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.011.10'
$ws.Range('E14').Value = '  +1.34%  '

# Row 15
$ws.Range('E15').Value = '  +0.52%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.611.98'
$ws.Range('E16').Value = '  +0.53%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.922'
$ws.Range('E17').Value = '  +1.71%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '15.01'
$ws.Range('E18').Value = '  -1.44%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '46.567.43'
$ws.Range('E19').Value = '  -1.40%  '

# Row 20
$ws.Range('E20').Value = '  +0.90%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.99'
$ws.Range('E21').Value = '  -8.19%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.75'
$ws.Range('E22').Value = '  +0.82%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.89'
$ws.Range('E23').Value = '  +1.74%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '275.14'
$ws.Range('E24').Value = '  +6.48%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.04'
$ws.Range('E25').Value = '  +1.19%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('E26').Value = '  +2.49%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '29.62'
$ws.Range('E27').Value = '  +15.19%  '

# Row 28
$ws.Range('E28').Value = '  -0.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.03'
$ws.Range('E29').Value = '  -1.30%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.62'
$ws.Range('E30').Value = '  +0.57%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '38.60'
$ws.Range('E31').Value = '  -8.81%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.22'
$ws.Range('E32').Value = '  -2.36%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.39'
$ws.Range('E33').Value = '  +5.76%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.65'
$ws.Range('E34').Value = '  -6.01%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.25'
$ws.Range('E35').Value = '  +0.24%  '

# Row 36
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.84'
$ws.Range('E36').Value = '  -4.18%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0840'
$ws.Range('E37').Value = '  -1.31%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.20'
$ws.Range('E38').Value = '  +1.11%  '

# Row 39
$ws.Range('E39').Value = '  +0.37%  '

# Row 40
$ws.Range('E40').Value = '  +1.13%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.39'
$ws.Range('E41').Value = '  +30.32%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.95'
$ws.Range('E42').Value = '  -4.17%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0332'
$ws.Range('E43').Value = '  +0.47%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.63'
$ws.Range('E44').Value = '  +0.09%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.04'
$ws.Range('E45').Value = '  -6.18%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.133.46'
$ws.Range('E46').Value = '  +6.06%  '

# Row 47
$ws.Range('E47').Value = '  +0.21%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '94.29'
$ws.Range('E48').Value = '  -1.03%  '

# Row 49
$ws.Range('E49').Value = '  +6.42%  '

# Row 50
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '109.95'
$ws.Range('E50').Value = '  +1.09%  '

# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.78'
$ws.Range('E51').Value = '  -7.93%  '
